$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.202022910118103
$ws.Range("B1").Value = 1.949655532836914
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.964297890663147
$ws.Range("E1").Value = 1.205873489379883
